# Scenario_EnergyCarrier_Price_CO2Emission.xlsx edit
# 1) Extend "Table2" on Sheet1 with 9 new year columns (2010..2018) inserted
#    right after the "unit" column and before the existing "2019" column,
#    shifting the existing 2019..2050 columns to the right.
# 2) Add a new worksheet "note" (after Sheet1) explaining the 2010-2018 values
#    are placeholders copied from 2019.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item(1)

function ColLetter([int]$n) {
    $letter = ""
    while ($n -gt 0) {
        $rem = ($n - 1) % 26
        $letter = [char](65 + $rem) + $letter
        $n = [int](($n - $rem - 1) / 26)
    }
    return $letter
}

# --- Step 1: grow the table by 9 columns (they land at the far right, AK:AS) ---
$newYearsCount = 9
for ($i = 0; $i -lt $newYearsCount; $i++) {
    $lo.ListColumns.Add() | Out-Null
}

# Existing year block occupies columns E..AJ (5..36), 32 columns (2019..2050).
# It needs to move to N..AS (14..45) to make room for the new 2010..2018 block.
$firstYearCol = 5          # E
$lastYearCol = 36           # AJ
$shift = $newYearsCount     # 9

# Copy cell-by-cell (value + number format) from right to left so we never
# overwrite a source cell before it has been read (dest = src + shift).
for ($c = $lastYearCol; $c -ge $firstYearCol; $c--) {
    $srcLetter = ColLetter $c
    $dstLetter = ColLetter ($c + $shift)
    for ($r = 1; $r -le 3; $r++) {
        $srcAddr = "$srcLetter$r"
        $dstAddr = "$dstLetter$r"
        $srcRange = $ws.Range($srcAddr)
        $val = $srcRange.Value()
        $fmt = $srcRange.NumberFormat
        $dstRange = $ws.Range($dstAddr)
        $dstRange.Value = $val
        $dstRange.NumberFormat = $fmt
    }
}

# --- Step 2: fill the freed-up E..M block with the new 2010..2018 data ---
# Header text "2010".."2018" ; body rows hold 0, formatted like the rest of
# the numeric data (style carried over from an existing data cell, e.g. N2).
$sampleFormat = $ws.Range("N2").NumberFormat
$years = 2010..2018
for ($i = 0; $i -lt $newYearsCount; $i++) {
    $c = $firstYearCol + $i
    $letter = ColLetter $c
    $ws.Range("$letter" + "1").Value = [string]$years[$i]
    for ($r = 2; $r -le 3; $r++) {
        $cell = $ws.Range("$letter$r")
        $cell.Value = 0
        $cell.NumberFormat = $sampleFormat
    }
}

# Give the new columns the same width as the rest of the sheet (10 chars).
$ws.Range("E1:M1").EntireColumn.ColumnWidth = 9.1

# Keep the table object in sync with the new right-hand edge.
$dims = $lo.Range.Address()
$lo.Resize($ws.Range("A1:AS3")) | Out-Null

# Restore the selection like the authored file (cosmetic, matches commit).
$ws.Range("D34").Select() | Out-Null

# --- Step 3: add the "note" worksheet after Sheet1 ---
$noteSheet = $wb.Worksheets.Add($null, $ws)
$noteSheet.Name = "note"
$noteSheet.Range("A1").Value = "note: the values from 2010 to 2018 are copied from 2019 to hold the place"
